# Replace unfinished menu items ("waffles") with a "placeholder" value in the
# Column1 column (F2:F6) of the pokebowls_greenmountain sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F6").Value = "placeholder"

# Leave the active selection on the last touched cell.
$ws.Range("F6").Select()
